$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- New shared strings (inserted via new cell content) ---
# Row 10: begin group / propriete_fonciere / II. PROPRIERE FONCIERE ET INFRASTRUCTURE DE BASE
$ws.Cells.Item(10, 1).Value = "begin group"
$ws.Cells.Item(10, 2).Value = "propriete_fonciere"
$ws.Cells.Item(10, 3).Value = "II. PROPRIERE FONCIERE ET INFRASTRUCTURE DE BASE"
$ws.Cells.Item(10, 4).Value = "II. PROPRIERE FONCIERE ET INFRASTRUCTURE DE BASE"

# Row 11: select_multiple source_elec / other_source_elec / fr / en
$ws.Cells.Item(11, 1).Value = "select_multiple source_elec"
$ws.Cells.Item(11, 2).Value = "other_source_elec"
$ws.Cells.Item(11, 3).Value = "Quelle(s) autre(s) source(s) d$([char]0x2019)électricité possède la FOSA ?"
$ws.Cells.Item(11, 4).Value = "What other source(s) of electricity does the health facility have?"

# Row 12: end group
$ws.Cells.Item(12, 1).Value = "end group"

# --- Formatting: apply "40% - Accent3" style look to B10:D10 ---
$style = $wb.Styles.Add("Excel Built-in 40% - Accent3")
$style.Font.Name = "Calibri"
$style.Font.Size = 11
$style.Font.Color = 0
$style.Interior.Color = 14149079
$ws.Range("B10:D10").Style = "Excel Built-in 40% - Accent3"

# Row heights to mirror the rest of the form's "13.8" rows, except the
# group-close row which uses the default "12.8".
$ws.Rows.Item(10).RowHeight = 13.8
$ws.Rows.Item(11).RowHeight = 13.8
$ws.Rows.Item(12).RowHeight = 12.8

# Re-use the existing "group row" / "label" styles for the new rows so
# they look the same as the other begin/end group + select_multiple rows.
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Style = "Normal"

# New column (C) used by the new group row's label columns.
$ws.Columns.Item(3).ColumnWidth = 25.22

# Selection cursor ends up on B18 after the edit (matches authoring session).
$ws.Range("B18").Select()

Write-Host "survey sheet updated"
